# FormatDataSiswa.xlsx - add "id_rombel" (rombel) values into column B
# for the data rows (rows 2-25), following the repeating 30/31/32 pattern,
# then restore the selection the author left active on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = 30, 31, 32

for ($r = 2; $r -le 25; $r++) {
    $idx = ($r - 2) % 3
    $ws.Cells.Item($r, 2).Value = $values[$idx]
}

# Re-select B5:B25 (with B5 as the active cell), matching the workbook's
# saved selection/view state after the edit.
$ws.Range("B5:B25").Select()
